$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "R-Homer_GARSTAR635_M-Bassoon_GAMAlexa594_0MgGlyBic_02.tif"
$ws.Range("B2").Value = 1900
$ws.Range("C2").Value = 994
$ws.Range("D2").Value = 321
$ws.Range("E2").Value = 0.05130075048611943
$ws.Range("F2").Value = 0.0980597846314154
$ws.Range("G2").Value = 8.131633999137502
